# cryptos.xlsx refresh - GitHub Actions scheduled update.
#
# Rewrites the Price (D) and Volume(1h) (E) columns for every coin row
# with freshly scraped figures, and also reflects a few rank swaps
# further down the table (rows 41-44: dogwifhat/Stacks and Maker/Fetch.AI
# traded places, bringing their Coin/Link/Price/Volume cells with them).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the A1 cell reference, the new literal text, and whether
# the text must be forced to stay text. Price strings such as "576.60"
# or "0.0000302" parse as plain numbers, and Excel would otherwise
# silently renormalize them (e.g. "576.60" -> 576.6, dropping the
# trailing zero) instead of keeping the scraped formatting verbatim.
$updates = @(
    @{ Cell = "D2"; Text = "70.193.02"; ForceText = $true }
    @{ Cell = "E2"; Text = "  -0.51%  "; ForceText = $false }
    @{ Cell = "D3"; Text = "3.589.22"; ForceText = $true }
    @{ Cell = "E3"; Text = "  -1.17%  "; ForceText = $false }
    @{ Cell = "E4"; Text = "  -0.03%  "; ForceText = $false }
    @{ Cell = "D5"; Text = "576.60"; ForceText = $true }
    @{ Cell = "E5"; Text = "  -2.97%  "; ForceText = $false }
    @{ Cell = "D6"; Text = "189.33"; ForceText = $true }
    @{ Cell = "E6"; Text = "  -1.80%  "; ForceText = $false }
    @{ Cell = "E7"; Text = "  -2.81%  "; ForceText = $false }
    @{ Cell = "D8"; Text = "3.586.05"; ForceText = $true }
    @{ Cell = "E8"; Text = "  -0.45%  "; ForceText = $false }
    @{ Cell = "E9"; Text = "  -0.01%  "; ForceText = $false }
    @{ Cell = "D10"; Text = "0.179"; ForceText = $true }
    @{ Cell = "E10"; Text = "  -1.45%  "; ForceText = $false }
    @{ Cell = "D11"; Text = "0.661"; ForceText = $true }
    @{ Cell = "E11"; Text = "  -0.38%  "; ForceText = $false }
    @{ Cell = "D12"; Text = "56.01"; ForceText = $true }
    @{ Cell = "E12"; Text = "  -3.45%  "; ForceText = $false }
    @{ Cell = "D13"; Text = "0.0000302"; ForceText = $true }
    @{ Cell = "E13"; Text = "  +1.10%  "; ForceText = $false }
    @{ Cell = "E14"; Text = "  -1.55%  "; ForceText = $false }
    @{ Cell = "D15"; Text = "4.162.86"; ForceText = $true }
    @{ Cell = "D16"; Text = "19.93"; ForceText = $true }
    @{ Cell = "E16"; Text = "  +2.80%  "; ForceText = $false }
    @{ Cell = "D17"; Text = "3.584.15"; ForceText = $true }
    @{ Cell = "E17"; Text = "  -1.36%  "; ForceText = $false }
    @{ Cell = "D18"; Text = "70.062.38"; ForceText = $true }
    @{ Cell = "E18"; Text = "  -0.49%  "; ForceText = $false }
    @{ Cell = "D19"; Text = "12.66"; ForceText = $true }
    @{ Cell = "E19"; Text = "  +0.22%  "; ForceText = $false }
    @{ Cell = "E20"; Text = "  +0.21%  "; ForceText = $false }
    @{ Cell = "E21"; Text = "  -1.06%  "; ForceText = $false }
    @{ Cell = "D22"; Text = "475.25"; ForceText = $true }
    @{ Cell = "E22"; Text = "  -4.44%  "; ForceText = $false }
    @{ Cell = "D23"; Text = "19.17"; ForceText = $true }
    @{ Cell = "E23"; Text = "  +14.51%  "; ForceText = $false }
    @{ Cell = "D24"; Text = "5.09"; ForceText = $true }
    @{ Cell = "E24"; Text = "  -7.96%  "; ForceText = $false }
    @{ Cell = "D25"; Text = "4.36"; ForceText = $true }
    @{ Cell = "E25"; Text = "  -2.04%  "; ForceText = $false }
    @{ Cell = "D26"; Text = "88.79"; ForceText = $true }
    @{ Cell = "E26"; Text = "  -2.40%  "; ForceText = $false }
    @{ Cell = "E27"; Text = "  -2.38%  "; ForceText = $false }
    @{ Cell = "D28"; Text = "11.06"; ForceText = $true }
    @{ Cell = "E28"; Text = "  -1.62%  "; ForceText = $false }
    @{ Cell = "D29"; Text = "9.37"; ForceText = $true }
    @{ Cell = "E29"; Text = "  -0.49%  "; ForceText = $false }
    @{ Cell = "D30"; Text = "32.19"; ForceText = $true }
    @{ Cell = "E30"; Text = "  -0.87%  "; ForceText = $false }
    @{ Cell = "E31"; Text = "  +2.09%  "; ForceText = $false }
    @{ Cell = "E32"; Text = "  +3.04%  "; ForceText = $false }
    @{ Cell = "D33"; Text = "12.15"; ForceText = $true }
    @{ Cell = "E33"; Text = "  -0.81%  "; ForceText = $false }
    @{ Cell = "D34"; Text = "66.18"; ForceText = $true }
    @{ Cell = "E34"; Text = "  +1.35%  "; ForceText = $false }
    @{ Cell = "D35"; Text = "584.32"; ForceText = $true }
    @{ Cell = "E35"; Text = "  -5.46%  "; ForceText = $false }
    @{ Cell = "E36"; Text = "  +2.63%  "; ForceText = $false }
    @{ Cell = "E37"; Text = "  -0.05%  "; ForceText = $false }
    @{ Cell = "D38"; Text = "0.0₃0799"; ForceText = $true }
    @{ Cell = "E38"; Text = "  -4.18%  "; ForceText = $false }
    @{ Cell = "D39"; Text = "0.398"; ForceText = $true }
    @{ Cell = "E39"; Text = "  -1.45%  "; ForceText = $false }
    @{ Cell = "E40"; Text = "  -6.41%  "; ForceText = $false }
    @{ Cell = "B41"; Text = "dogwifhat"; ForceText = $false }
    @{ Cell = "C41"; Text = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; ForceText = $false }
    @{ Cell = "D41"; Text = "3.22"; ForceText = $true }
    @{ Cell = "E41"; Text = "  +16.50%  "; ForceText = $false }
    @{ Cell = "B42"; Text = "Stacks"; ForceText = $false }
    @{ Cell = "C42"; Text = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; ForceText = $false }
    @{ Cell = "D42"; Text = "3.48"; ForceText = $true }
    @{ Cell = "E42"; Text = "  -6.12%  "; ForceText = $false }
    @{ Cell = "B43"; Text = "Maker"; ForceText = $false }
    @{ Cell = "C43"; Text = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; ForceText = $false }
    @{ Cell = "D43"; Text = "3.244.06"; ForceText = $true }
    @{ Cell = "E43"; Text = "  -3.31%  "; ForceText = $false }
    @{ Cell = "B44"; Text = "Fetch.AI"; ForceText = $false }
    @{ Cell = "C44"; Text = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; ForceText = $false }
    @{ Cell = "D44"; Text = "2.87"; ForceText = $true }
    @{ Cell = "E44"; Text = "  +7.27%  "; ForceText = $false }
    @{ Cell = "E45"; Text = "  +1.33%  "; ForceText = $false }
    @{ Cell = "D46"; Text = "0.0441"; ForceText = $true }
    @{ Cell = "E46"; Text = "  -1.38%  "; ForceText = $false }
    @{ Cell = "D47"; Text = "9.51"; ForceText = $true }
    @{ Cell = "E47"; Text = "  +4.31%  "; ForceText = $false }
    @{ Cell = "E48"; Text = "  -0.49%  "; ForceText = $false }
    @{ Cell = "E49"; Text = "  -0.63%  "; ForceText = $false }
    @{ Cell = "D50"; Text = "0.998"; ForceText = $true }
    @{ Cell = "E50"; Text = "  -0.18%  "; ForceText = $false }
    @{ Cell = "D51"; Text = "3.15"; ForceText = $true }
    @{ Cell = "E51"; Text = "  -3.67%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe is Excel's own "treat as text" quote
        # prefix; it is not part of the stored value.
        $cell.Value = "'" + $u.Text
    } else {
        $cell.Value = $u.Text
    }
}
